$d = $word.ActiveDocument

# --- Add the three new character styles ---

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply GaNStyle to every "Kaksosten tähtikuvio havainnointijaksot..." run,
#     also fixing the missing trailing period ---

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "Kaksosten tähtikuvio havainnointijaksot vuonna 2022: 14.-23.2., 14.-24.3"
$rng.Find.Forward = $true
$rng.Find.Wrap = 1
while ($rng.Find.Execute()) {
    $rng.Text = "Kaksosten tähtikuvio havainnointijaksot vuonna 2022: 14.-23.2., 14.-24.3."
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- Apply GaNParagraph to the "Osallistut maailmanlaajuiseen kampanjaan..." run ---

$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Text = "Osallistut maailmanlaajuiseen kampanjaan tarkkaillaksesi ja tallentaaksesi himmeimpiä näkyvissä olevia tähtiä keinona mitata valonsaastetta tietyssä paikassa. Paikallistamalla ja tarkkailemalla Kaksosten tähtikuvio miten valosaaste syntyy kunkin taajaman tai muun ihmisen toiminnan valoista. Antamasi tiedot päivittyvät heti verkossa olevaan tietokantaan, ja näin saadaan käsitys siitä minkä verran taivaan tähdistä on missäkin nähtävissä."
$rng2.Find.Forward = $true
$rng2.Find.Wrap = 1
if ($rng2.Find.Execute()) {
    $rng2.Style = "GaNParagraph"
}

# --- Apply GaNLinks to the "Tämän oppaan kartat piirsi..." run ---

$rng3 = $d.Content
$rng3.Find.ClearFormatting()
$rng3.Find.Text = "Tämän oppaan kartat piirsi Jenik Hollan CzechGlobesta (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."
$rng3.Find.Forward = $true
$rng3.Find.Wrap = 1
if ($rng3.Find.Execute()) {
    $rng3.Style = "GaNLinks"
}
